$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings are not
# auto-converted to numbers, matching the original inlineStr (text) cell content.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = '29.430.77'
$ws.Range("E2").Value = '  -0.22%  '

$ws.Range("D3").Value = '1.848.28'
$ws.Range("E3").Value = '  -0.23%  '

$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("E5").Value = '  -1.03%  '

$ws.Range("D6").Value = '0.6337'
$ws.Range("E6").Value = '  -0.24%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Value = '3.653.95'
$ws.Range("E8").Value = '  +97.14%  '

$ws.Range("D9").Value = '0.07555'
$ws.Range("E9").Value = '  +1.14%  '

$ws.Range("D10").Value = '0.2971'
$ws.Range("E10").Value = '  -0.79%  '

$ws.Range("B11").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C11").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D11").Value = '3.915.38'
$ws.Range("E11").Value = '  +87.32%  '

$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").Value = '24.64'
$ws.Range("E12").Value = '  +1.55%  '

$ws.Range("D13").Value = '0.07711'
$ws.Range("E13").Value = '  +1.04%  '

$ws.Range("D14").Value = '4.990'
$ws.Range("E14").Value = '  -0.74%  '

$ws.Range("D15").Value = '0.6859'
$ws.Range("E15").Value = '  +0.05%  '

$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = '82.99'
$ws.Range("E16").Value = '  -0.69%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.000009934'
$ws.Range("E17").Value = '  +4.25%  '

$ws.Range("D18").Value = '6.188'
$ws.Range("E18").Value = '  +0.36%  '

$ws.Range("D19").Value = '29.468.26'
$ws.Range("E19").Value = '  -0.13%  '

$ws.Range("D20").Value = '231.71'
$ws.Range("E20").Value = '  -1.58%  '

$ws.Range("D21").Value = '12.49'
$ws.Range("E21").Value = '  -0.43%  '

$ws.Range("D22").Value = '0.9998'
$ws.Range("E22").Value = '  -0.03%  '

$ws.Range("D23").Value = '7.597'
$ws.Range("E23").Value = '  -1.03%  '

$ws.Range("E24").Value = '  -0.07%  '

$ws.Range("D25").Value = '155.51'

$ws.Range("E26").Value = '  -1.24%  '

$ws.Range("D27").Value = '8.425'
$ws.Range("E27").Value = '  -0.81%  '

$ws.Range("E28").Value = '  -0.38%  '

$ws.Range("B29").Value = 'RocketPoolETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D29").Value = '3.893.37'
$ws.Range("E29").Value = '  +94.00%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '1.467'
$ws.Range("E30").Value = '  -1.38%  '

$ws.Range("D31").Value = '0.05817'
$ws.Range("E31").Value = '  -3.19%  '

$ws.Range("E32").Value = '  +0.69%  '

$ws.Range("D33").Value = '4.140'
$ws.Range("E33").Value = '  +0.38%  '

$ws.Range("D34").Value = '4.023'
$ws.Range("E34").Value = '  -1.30%  '

$ws.Range("D35").Value = '1.862'
$ws.Range("E35").Value = '  -0.21%  '

$ws.Range("D36").Value = '1.156'
$ws.Range("E36").Value = '  -1.99%  '

$ws.Range("D37").Value = '0.7166'
$ws.Range("E37").Value = '  -0.29%  '

$ws.Range("D38").Value = '2.592'
$ws.Range("E38").Value = '  -0.21%  '

$ws.Range("D39").Value = '1.252.62'
$ws.Range("E39").Value = '  +4.30%  '

$ws.Range("D40").Value = '2.797'
$ws.Range("E40").Value = '  -0.11%  '

$ws.Range("D41").Value = '0.01805'
$ws.Range("E41").Value = '  +1.80%  '

$ws.Range("D42").Value = '0.9013'
$ws.Range("E42").Value = '  -0.71%  '

$ws.Range("D43").Value = '6.100'
$ws.Range("E43").Value = '  -1.06%  '

$ws.Range("D44").Value = '0.9998'
$ws.Range("E44").Value = '  +0.03%  '

$ws.Range("D45").Value = '101.79'

$ws.Range("D46").Value = '67.05'
$ws.Range("E46").Value = '  +0.81%  '

$ws.Range("D47").Value = '7.210'
$ws.Range("E47").Value = '  -1.29%  '

$ws.Range("D48").Value = '9.144'
$ws.Range("E48").Value = '  +0.54%  '

$ws.Range("D49").Value = '0.4017'
$ws.Range("E49").Value = '  -0.35%  '

$ws.Range("D50").Value = '1.685'
$ws.Range("E50").Value = '  +1.62%  '

$ws.Range("E51").Value = '  +0.15%  '

# Restore default style (remove the temporary text-format override) so
# the cells retain their original (unstyled) appearance.
$rng.Style = "Normal"
